$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 10000
$ws.Range("J47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("N47").Value = -11944
$ws.Range("H68").Value = 35295
$ws.Range("J68").Value = 35295
$ws.Range("L68").Value = 35295
$ws.Range("N68").Value = -36793
$ws.Range("H71").Value = 35295
$ws.Range("J71").Value = 35295
$ws.Range("L71").Value = 105885
$ws.Range("N71").Value = -113373

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 50003204
$ws.Range("I2").Value = 83334340
$ws.Range("J2").Value = 6506.5
$ws.Range("K2").Value = 83334340
$ws.Range("L2").Value = 6506.5
$ws.Range("M2").Value = -83334227
$ws.Range("N2").Value = -6732.5
$ws.Range("H32").Value = 4653.837
$ws.Range("I32").Value = 3567.7927
$ws.Range("K32").Value = 3567.7927
$ws.Range("M32").Value = -3280.7927
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H45").Value = 1558.7693
$ws.Range("I45").Value = 996.5161000000001
$ws.Range("J45").Value = 3737.5
$ws.Range("K45").Value = 996.5161000000001
$ws.Range("L45").Value = 3737.5
$ws.Range("M45").Value = -619.5161000000001
$ws.Range("N45").Value = -4491.5
$ws.Range("H97").Value = 615.5769
$ws.Range("I97").Value = 582.8261
$ws.Range("K97").Value = 582.8261
$ws.Range("M97").Value = -86.8261
$ws.Range("H116").Value = 50003204
$ws.Range("I116").Value = 83334340
$ws.Range("J116").Value = 6506.5
$ws.Range("K116").Value = 83334340
$ws.Range("L116").Value = 6506.5
$ws.Range("M116").Value = -83332046
$ws.Range("N116").Value = -11094.5
$ws.Range("H132").Value = 22224874
$ws.Range("I132").Value = 33335070
$ws.Range("J132").Value = 4483.467
$ws.Range("K132").Value = 100005210
$ws.Range("L132").Value = 13450.401
$ws.Range("M132").Value = -100002680
$ws.Range("N132").Value = -18510.401
$ws.Range("H139").Value = 29562.857
$ws.Range("J139").Value = 29562.857
$ws.Range("L139").Value = 29562.857
$ws.Range("N139").Value = -39842.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 50003204
$ws.Range("I3").Value = 83334340
$ws.Range("J3").Value = 6506.5
$ws.Range("K3").Value = 83334340
$ws.Range("L3").Value = 6506.5
$ws.Range("M3").Value = -83334226
$ws.Range("N3").Value = -6734.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2564.5745
$ws.Range("I31").Value = 1988.6586
$ws.Range("J31").Value = 6500
$ws.Range("K31").Value = 1988.6586
$ws.Range("L31").Value = 6500
$ws.Range("M31").Value = -1693.6586
$ws.Range("N31").Value = -7090
$ws.Range("H34").Value = 2564.5745
$ws.Range("I34").Value = 1988.6586
$ws.Range("J34").Value = 6500
$ws.Range("K34").Value = 1988.6586
$ws.Range("L34").Value = 6500
$ws.Range("M34").Value = -1786.6586
$ws.Range("N34").Value = -6904
$ws.Range("H59").Value = 7675.2856
$ws.Range("J59").Value = 8287.833000000001
$ws.Range("L59").Value = 8287.833000000001
$ws.Range("N59").Value = -10577.833
$ws.Range("H62").Value = 4278
$ws.Range("I62").Value = 2526.6667
$ws.Range("J62").Value = 5028.5713
$ws.Range("K62").Value = 2526.6667
$ws.Range("L62").Value = 5028.5713
$ws.Range("M62").Value = -1902.6667
$ws.Range("N62").Value = -6276.5713
$ws.Range("H65").Value = 4278
$ws.Range("I65").Value = 2526.6667
$ws.Range("J65").Value = 5028.5713
$ws.Range("K65").Value = 12633.3335
$ws.Range("L65").Value = 25142.8565
$ws.Range("M65").Value = -9513.333500000001
$ws.Range("N65").Value = -31382.8565
$ws.Range("H132").Value = 2882.963
$ws.Range("I132").Value = 2207.3684
$ws.Range("J132").Value = 4487.5
$ws.Range("K132").Value = 6622.1052
$ws.Range("L132").Value = 13462.5
$ws.Range("M132").Value = -4092.1052
$ws.Range("N132").Value = -18522.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1879.625
$ws.Range("I59").Value = 284.25
$ws.Range("J59").Value = 3475
$ws.Range("K59").Value = 852.75
$ws.Range("L59").Value = 10425
$ws.Range("M59").Value = -312.75
$ws.Range("N59").Value = -11505
$ws.Range("H87").Value = 13425
$ws.Range("J87").Value = 15677.777
$ws.Range("L87").Value = 47033.331
$ws.Range("N87").Value = -49529.331
$ws.Range("H90").Value = 13425
$ws.Range("J90").Value = 15677.777
$ws.Range("L90").Value = 141099.993
$ws.Range("N90").Value = -153579.993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3746.0715
$ws.Range("I80").Value = 4016.5
$ws.Range("J80").Value = 3070
$ws.Range("K80").Value = 4016.5
$ws.Range("L80").Value = 3070
$ws.Range("M80").Value = -3018.5
$ws.Range("N80").Value = -5066
$ws.Range("H83").Value = 3746.0715
$ws.Range("I83").Value = 4016.5
$ws.Range("J83").Value = 3070
$ws.Range("K83").Value = 20082.5
$ws.Range("L83").Value = 15350
$ws.Range("M83").Value = -15090.5
$ws.Range("N83").Value = -25334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1480
$ws.Range("I68").Value = 1023.8095
$ws.Range("J68").Value = 3875
$ws.Range("K68").Value = 1023.8095
$ws.Range("L68").Value = 3875
$ws.Range("M68").Value = -274.8095
$ws.Range("N68").Value = -5373
$ws.Range("H71").Value = 1480
$ws.Range("I71").Value = 1023.8095
$ws.Range("J71").Value = 3875
$ws.Range("K71").Value = 5119.0475
$ws.Range("L71").Value = 19375
$ws.Range("M71").Value = -1375.0475
$ws.Range("N71").Value = -26863
$ws.Range("H82").Value = 3125
$ws.Range("I82").Value = 2142.8572
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 2142.8572
$ws.Range("L82").Value = 10000
$ws.Range("M82").Value = -1781.8572
$ws.Range("N82").Value = -10722
$ws.Range("H85").Value = 3125
$ws.Range("I85").Value = 2142.8572
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 2142.8572
$ws.Range("L85").Value = 10000
$ws.Range("M85").Value = -894.8571999999999
$ws.Range("N85").Value = -12496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 35097.145
$ws.Range("J16").Value = 35097.145
$ws.Range("L16").Value = 35097.145
$ws.Range("N16").Value = -35681.145
